$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = "2/1"
$ws.Range("F3").Value = "2/0"
$ws.Range("G3").Value = "1/2"
